$wb = $excel.ActiveWorkbook

# Fix typo in sheet name: "Dimensios" -> "Dimensions"
$dimSheet = $wb.Worksheets.Item("Dimensios")
$dimSheet.Name = "Dimensions"

# Saboteurs sheet: update D2 value/style (match C2's yellow fill + border) and selection
$sabSheet = $wb.Worksheets.Item("Saboteurs")
$sabSheet.Activate()
$d2 = $sabSheet.Range("D2")
$d2.Value = 2
$d2.Interior.Pattern = -4105
$d2.Interior.Color = 65535
$sabSheet.Range("C2:D2").Select()

# Switch the active tab from "Dimensions" to "Obstacles"
$obsSheet = $wb.Worksheets.Item("Obstacles")
$obsSheet.Activate()
